# Update gh-pages to output generated at 456a3b4
# Updates the F column ('想去人数' / want-to-go counts) across the four sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 546
$ws.Range("F5").Value = 82
$ws.Range("F6").Value = 3705
$ws.Range("F11").Value = 1541
$ws.Range("F12").Value = 8
$ws.Range("F14").Value = 641
$ws.Range("F15").Value = 1529
$ws.Range("F16").Value = 1405
$ws.Range("F17").Value = 25
$ws.Range("F18").Value = 28
$ws.Range("F19").Value = 556
$ws.Range("F20").Value = 4011
$ws.Range("F21").Value = 4012
$ws.Range("F22").Value = 660
$ws.Range("F23").Value = 3310
$ws.Range("F24").Value = 765
$ws.Range("F25").Value = 34
$ws.Range("F26").Value = 2222
$ws.Range("F28").Value = 314
$ws.Range("F30").Value = 33
$ws.Range("F31").Value = 1176
$ws.Range("F33").Value = 50
$ws.Range("F34").Value = 1067
$ws.Range("F35").Value = 1075

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 46
$ws.Range("F6").Value = 96
$ws.Range("F13").Value = 213
$ws.Range("F18").Value = 264
$ws.Range("F19").Value = 198

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 243
$ws.Range("F4").Value = 547
$ws.Range("F5").Value = 73
$ws.Range("F6").Value = 165

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 243
$ws.Range("F4").Value = 546
$ws.Range("F5").Value = 46
$ws.Range("F8").Value = 82
$ws.Range("F9").Value = 547
$ws.Range("F10").Value = 3708
$ws.Range("F11").Value = 3709
$ws.Range("F17").Value = 96
$ws.Range("F22").Value = 1541
$ws.Range("F25").Value = 1529
$ws.Range("F27").Value = 1405
$ws.Range("F28").Value = 25
$ws.Range("F29").Value = 556
$ws.Range("F31").Value = 4012
$ws.Range("F32").Value = 4012
$ws.Range("F33").Value = 660
$ws.Range("F34").Value = 3310
$ws.Range("F35").Value = 765
$ws.Range("F36").Value = 34
$ws.Range("F37").Value = 2222
$ws.Range("F39").Value = 314
$ws.Range("F41").Value = 33
$ws.Range("F42").Value = 1176
$ws.Range("F44").Value = 264
$ws.Range("F45").Value = 198
$ws.Range("F48").Value = 50
$ws.Range("F49").Value = 1067
$ws.Range("F50").Value = 1075
